$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 2 - this pushes the existing
# classifier rows (originally rows 2-13) down to rows 3-14, carrying
# each row's own formatting along with it.
$ws.Rows.Item(2).Insert() | Out-Null

# Clone formatting onto the new row 2 before filling it in:
#  - B2:C2 take on the numeric/left-aligned style used throughout column B/C
#    (copy it from B3, which is the old row-2 cell that just got pushed down).
#  - A2 takes on the plain/default style used by the other classifier-name
#    cells (copy it from A4, an existing unstyled cell).
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B2:C2").PasteSpecial(-4122) | Out-Null
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null

# Fill in the new GradientBoostingClassifier row.
$ws.Range("A2").Value = "GradientBoostingClassifier"
$ws.Range("B2").Value = 0.80342240975152301
$ws.Range("C2").Value = "{'clf__n_estimators': 500}"

# Re-sort the full data range (A2:C14) descending by Accuracy (col B),
# same as the author re-sorting the table after adding the new row.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear() | Out-Null
$sortObj.SortFields.Add($ws.Range("B2:B14"), $null, 2, $null, 0) | Out-Null
$sortObj.SetRange($ws.Range("A2:C14")) | Out-Null
$sortObj.Header = 0
$sortObj.Apply() | Out-Null

# GaussianNB (now row 14) has no Hyper Parameters value - drop the blank
# styled cell the sort leaves behind in C14 so it matches a truly empty cell.
$ws.Range("C14").Clear() | Out-Null

# Leave the cursor where the author left it when they saved the file.
$ws.Range("C9").Select() | Out-Null
